$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 495, shifting rows 495:519 down to 496:520
$ws.Rows.Item(495).Insert()

# Fill in the new row 495 with the same categorical values as its neighboring
# carrot ("Zanahoria") records, and the new date/volume/price data point.
$ws.Cells.Item(495, 1).Value = 8
$ws.Cells.Item(495, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(495, 3).Value = "Coquimbo"
$ws.Cells.Item(495, 4).Value = 45041
$ws.Cells.Item(495, 4).NumberFormat = $ws.Cells.Item(496, 4).NumberFormat
$ws.Cells.Item(495, 5).Value = 4
$ws.Cells.Item(495, 6).Value = 100114013
$ws.Cells.Item(495, 7).Value = "Zanahoria"
$ws.Cells.Item(495, 8).Value = "Sin especificar"
$ws.Cells.Item(495, 9).Value = "Primera"
$ws.Cells.Item(495, 10).Value = 600
$ws.Cells.Item(495, 11).Value = 5500
$ws.Cells.Item(495, 12).Value = 6000
$ws.Cells.Item(495, 13).Value = 5750
$ws.Cells.Item(495, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(495, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(495, 16).Value = 288
$ws.Cells.Item(495, 17).Value = 20
$ws.Cells.Item(495, 18).Value = "Hortaliza"
